$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "35÷2=17, 1"
$t.Cell(1, 2).Range.Text = "72÷9=8, 0"
$t.Cell(1, 3).Range.Text = "37÷5=7, 2"
$t.Cell(1, 4).Range.Text = "51÷3=17, 0"
$t.Cell(1, 5).Range.Text = "52÷9=5, 7"
$t.Cell(5, 1).Range.Text = "59÷6=9, 5"
$t.Cell(5, 2).Range.Text = "28÷5=5, 3"
$t.Cell(5, 3).Range.Text = "49÷9=5, 4"
$t.Cell(5, 4).Range.Text = "50÷5=10, 0"
$t.Cell(5, 5).Range.Text = "94÷9=10, 4"
$t.Cell(9, 1).Range.Text = "90÷6=15, 0"
$t.Cell(9, 2).Range.Text = "82÷8=10, 2"
$t.Cell(9, 3).Range.Text = "95÷6=15, 5"
$t.Cell(9, 4).Range.Text = "25÷2=12, 1"
$t.Cell(9, 5).Range.Text = "93÷3=31, 0"
$t.Cell(13, 1).Range.Text = "49÷5=9, 4"
$t.Cell(13, 2).Range.Text = "56÷4=14, 0"
$t.Cell(13, 3).Range.Text = "12÷6=2, 0"
$t.Cell(13, 4).Range.Text = "59÷7=8, 3"
$t.Cell(13, 5).Range.Text = "40÷9=4, 4"
$t.Cell(17, 1).Range.Text = "71÷6=11, 5"
$t.Cell(17, 2).Range.Text = "55÷8=6, 7"
$t.Cell(17, 3).Range.Text = "25÷2=12, 1"
$t.Cell(17, 4).Range.Text = "72÷9=8, 0"
$t.Cell(17, 5).Range.Text = "89÷7=12, 5"
